$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Reverse Linked List
$ws.Range("A4").Value = 206
$ws.Range("B4").Value = "Reverse Linked List"
$ws.Range("C4").Value = "Linked List"
$ws.Range("D4").Value = "Pedoe"
$ws.Range("E4").Value = "Ongoing"
$ws.Range("F4").Value = "Easy"
$ws.Range("G4").Value = "Reverse a single linked list"

# Row 5: Linked List Cycle
$ws.Range("A5").Value = 141
$ws.Range("B5").Value = "Linked List Cycle"
$ws.Range("C5").Value = "Linked List"
$ws.Range("D5").Value = "Pedoe"
$ws.Range("E5").Value = "Ongoing"
$ws.Range("F5").Value = "Easy"
$ws.Range("G5").Value = "Solve it without using extra space"

# Row 6: Merge Two Sort Lists
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = "Merge Two Sort Lists"
$ws.Range("C6").Value = "Linked List"
$ws.Range("D6").Value = "Pedoe"
$ws.Range("E6").Value = "Ongoing"
$ws.Range("F6").Value = "Easy"

# Copy the "Ongoing" style from E3 onto E4:E6 (same cellXf reused, no new style)
$ws.Range("E3").Copy()
$ws.Range("E4:E6").PasteSpecial(-4122)

$ws.Range("E11").Select()
